$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:AXISBANK"
$ws.Range("C2").Value = "NSE:AARVI"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "NSE:HAL"
$ws.Range("F2").Value = "NSE:BRITANNIA"

# Row 3
$ws.Range("B3").Value = "NSE:BRITANNIA"
$ws.Range("C3").Value = "NSE:ADANIENT"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

# Row 4
$ws.Range("B4").Value = "NSE:COLPAL"
$ws.Range("C4").Value = "NSE:AGI"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("B5").Value = "NSE:CONCORDBIO"
$ws.Range("C5").Value = "NSE:ALKYLAMINE"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("B6").Value = "NSE:HEXATRADEX"
$ws.Range("C6").Value = "NSE:AURUM"

# Row 7
$ws.Range("B7").Value = "NSE:LUXIND"
$ws.Range("C7").Value = "NSE:BANSWRAS"

# Row 8
$ws.Range("B8").Value = "NSE:MANAKALUCO"
$ws.Range("C8").Value = "NSE:BEPL"

# Row 9
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "NSE:CENTENKA"

# Row 10
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "NSE:COCHINSHIP"

# Row 11
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "NSE:COFORGE"

# Row 12
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "NSE:DCAL"

# Row 13
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "NSE:DREDGECORP"

# Row 14
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "NSE:DYNAMATECH"

# Row 15
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:EVERESTIND"

# Row 16
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "NSE:FAIRCHEMOR"

# Row 17
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "NSE:GALAXYSURF"

# Row 18
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "NSE:GRSE"

# Row 19
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = "NSE:GRWRHITECH"

# Row 20
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = "NSE:GTLINFRA"

# Row 21
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = "NSE:HINDZINC"

# Row 22
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "NSE:IFBIND"

# Row 23
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = "NSE:KAMATHOTEL"

# Row 24
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "NSE:KBCGLOBAL"

# Row 25
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = "NSE:KRISHANA"

# Row 26
$ws.Range("B26").Value = ""
$ws.Range("C26").Value = "NSE:MAZDOCK"

# Row 27
$ws.Range("B27").Value = ""
$ws.Range("C27").Value = "NSE:MIRZAINT"

# Row 28
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = "NSE:MTARTECH"

# Row 29
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = "NSE:NRAIL"

# Row 30
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = "NSE:ONEPOINT"

# Row 31
$ws.Range("B31").Value = ""
$ws.Range("C31").Value = "NSE:PRECOT"

# Row 32
$ws.Range("B32").Value = ""
$ws.Range("C32").Value = "NSE:PRITIKAUTO"

# Row 33
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = "NSE:RKFORGE"

# Remove row 34 entirely (was NSE:RADIANTCMS), shifting dimension to A1:F33
$ws.Range("A34:F34").Delete()
